$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39
$src = 38

# Values first
$ws.Cells.Item($row, 1).Value = 38
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45257.64583333334
$ws.Cells.Item($row, 6).Value = "Goa"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Jamshedpur"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 1.51
$ws.Cells.Item($row, 11).Value = "20/11/2023 15:42"
$ws.Cells.Item($row, 12).Value = 1.52
$ws.Cells.Item($row, 13).Value = "27/11/2023 15:26"
$ws.Cells.Item($row, 14).Value = 4.39
$ws.Cells.Item($row, 15).Value = "20/11/2023 15:42"
$ws.Cells.Item($row, 16).Value = 4.54
$ws.Cells.Item($row, 17).Value = "27/11/2023 15:26"
$ws.Cells.Item($row, 18).Value = 5.82
$ws.Cells.Item($row, 19).Value = "20/11/2023 15:42"
$ws.Cells.Item($row, 20).Value = 5.75
$ws.Cells.Item($row, 21).Value = "27/11/2023 15:26"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/fc-goa-jamshedpur/ljShS5h2/"

# Copy formatting from the row above (A column is bold/bordered/centered,
# E column carries the datetime number format) so no new styles are minted.
$ws.Cells.Item($src, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($src, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
